# Coba import data penjualan
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the penjualan_kode (B) for both rows first, then pembeli (C), then barang_id (D)
$ws.Range("B2").Value = "TXR0022"
$ws.Range("B3").Value = "TXR0023"
$ws.Range("C2").Value = "Customer 22"
$ws.Range("C3").Value = "Customer 23"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1

# Move selection to E4 (matches the saved cursor position in the diff)
$ws.Range("E4").Select()
